# Fruta / hortaliza, semanal
# Insert 3 new price rows (a new reporting date) above the existing
# row 300 in the "Femacal de La Calera - Chirimoya" data table, pushing
# the rest of the table down by 3 rows (old 300..312 -> 303..315).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before row 300.
$ws.Rows.Item(300).Insert()
$ws.Rows.Item(300).Insert()
$ws.Rows.Item(300).Insert()

# Common template values shared by every row of this data block.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$tipo      = "Fruta"
$productoId = 100107
$producto   = "Otros"
$categoriaId = 100107002
$categoria   = "Chirimoya"
$variedad    = "Cultivar IV Región"
$unidad      = "$/bandeja 10 kilos"
$kgUnidad    = 10

$fecha = 45147
$origen = "Provincia del Elquí"

# Row 300 - Especial
$r = 300
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = $ws.Range("D303").NumberFormat
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 40
$ws.Cells.Item($r, 14).Value = 33000
$ws.Cells.Item($r, 15).Value = 33000
$ws.Cells.Item($r, 16).Value = 33000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 3300
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 301 - Primera
$r = 301
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = $ws.Range("D304").NumberFormat
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 35
$ws.Cells.Item($r, 14).Value = 30000
$ws.Cells.Item($r, 15).Value = 30000
$ws.Cells.Item($r, 16).Value = 30000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 3000
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 302 - Segunda
$r = 302
$ws.Cells.Item($r, 1).Value = $mercadoId
$ws.Cells.Item($r, 2).Value = $mercado
$ws.Cells.Item($r, 3).Value = $region
$ws.Cells.Item($r, 4).Value = $fecha
$ws.Cells.Item($r, 4).NumberFormat = $ws.Range("D305").NumberFormat
$ws.Cells.Item($r, 5).Value = $codreg
$ws.Cells.Item($r, 6).Value = $tipo
$ws.Cells.Item($r, 7).Value = $productoId
$ws.Cells.Item($r, 8).Value = $producto
$ws.Cells.Item($r, 9).Value = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = $variedad
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 36
$ws.Cells.Item($r, 14).Value = 28000
$ws.Cells.Item($r, 15).Value = 28000
$ws.Cells.Item($r, 16).Value = 28000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 2800
$ws.Cells.Item($r, 20).Value = $kgUnidad
